$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps literal text values (no numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.742.34"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.696.40"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "316.54"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.3952"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "1.488"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "1.000"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "53.23"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "0.08858"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "7.275"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "23.70"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "8.060"
$ws.Range("E15").Value = "  +8.62%  "
$ws.Range("D16").Value = "0.00001321"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "1.696.14"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "100.11"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "0.07026"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "19.59"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "7.052"
$ws.Range("E21").Value = "  +4.82%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "14.34"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "24.736.47"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "3.257"
$ws.Range("E25").Value = "  +9.46%  "
$ws.Range("D26").Value = "2.359"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").Value = "22.85"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "163.76"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").Value = "136.34"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "5.182"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").Value = "7.520"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "1.883.83"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "1.073"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "0.08601"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "7.139"
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("D36").Value = "11.43"
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("D37").Value = "0.2750"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").Value = "1.918"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "14.49"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "0.09225"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").Value = "0.02731"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "1.467"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "0.7683"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").Value = "16.12"
$ws.Range("E44").Value = "  +4.58%  "
$ws.Range("D45").Value = "0.7196"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "2.589"
$ws.Range("E46").Value = "  +6.40%  "
$ws.Range("D47").Value = "4.217"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "1.325"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "139.78"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "0.07986"
$ws.Range("E51").Value = "  +0.96%  "

Write-Output "Updated cryptos list"
